# updated legacy GSC export data
# The daily GSC export rolled forward by one day: the oldest day
# (2025-11-01) drops off the top of the table and three new days
# (2026-01-27, 2026-01-28, 2026-01-29) are appended at the bottom.
# Deleting the first data row shifts every remaining row up by one,
# which is exactly equivalent to the whole date range advancing by a
# day while each row keeps the metric values that used to belong to
# the row below it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Drop the oldest day (row 2 = 2025-11-01); everything below shifts up.
$ws.Rows(2).Delete()

# Helper: write a literal text date (not an auto-converted date serial)
# into a cell without leaving behind a new cell style. We do this by
# computing the text via a TEXT() formula, then flattening the formula
# to its literal value with a values-only paste.
function Set-TextDate($cell, [int]$y, [int]$m, [int]$d) {
    $formula = '=TEXT(DATE(' + $y + ',' + $m + ',' + $d + '),"yyyy-MM-dd")'
    $cell.Formula = $formula
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}

# Row 88: newly observed day 2026-01-27
Set-TextDate $ws.Range("A88") 2026 1 27
$ws.Range("B88").Value = 0
$ws.Range("C88").Value = 0
$ws.Range("D88").Value = 0

# Row 89: newly observed day 2026-01-28
Set-TextDate $ws.Range("A89") 2026 1 28
$ws.Range("B89").Value = 0
$ws.Range("C89").Value = 0
$ws.Range("D89").Value = 0

# Row 90: newly observed day 2026-01-29 (impressions not yet populated,
# matching the same blank placeholder used for other very recent days)
Set-TextDate $ws.Range("A90") 2026 1 29
$ws.Range("B90").Value = 0
$ws.Range("C90").Value = 0
$ws.Range("D90").Value = ""
